# Regenerate merged AHB files
# - Rename the "_old"/"_new" suffixed header labels (row 1, columns A:J and
#   L:U) to "_FV2304"/"_FV2310" respectively.
# - Turn the used range A1:U93 into an Excel Table ("Table1").
# - Freeze the header row (row 1) in the worksheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (A1:J1 "_old" -> "_FV2304", L1:U1 "_new" -> "_FV2310") ---
$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"

# K1 ("diff") is left untouched.

$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# --- Convert A1:U93 into a native Excel table (adds xl/tables/table1.xml,
#     the sheet -> table relationship, tableParts entry, and the Content
#     Types override) ---
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U93"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- Freeze the header row (pane split after row 1) ---
$null = $ws.Activate()
$null = $ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
